$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed value
$ws.Range("A7").Value = 1

# Diagonal chain of non-shared formulas
$ws.Range("B8").Formula = "=A7+1"
$ws.Range("C9").Formula = "=B8+1"
$ws.Range("D10").Formula = "=C9+1"
$ws.Range("E11").Formula = "=D10+1"
$ws.Range("F12").Formula = "=E11+1"
$ws.Range("G13").Formula = "=F12+1"
$ws.Range("H14").Formula = "=G13+1"
$ws.Range("I15").Formula = "=H14+1"
$ws.Range("J16").Formula = "=I15+1"

# Select the last active cell as in the final sheet state
$ws.Range("J16").Select()
